$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting of the new rows (13..31) in column A to match the
# existing numbered rows (style index used by A2:A12 - bold, bordered,
# centered/top aligned).
$ws.Range("A12").Copy()
$ws.Range("A13:A31").PasteSpecial(-4122)

# Full target data for rows 2..31 across columns A..F.
# Column A: sequential index 0..29
# Column B: only populated for rows 2..5
# Column C: alphabetically sorted ticker list (now 25 entries, rows 2..26... actually 2..31 minus blanks)
# Column D: always blank
# Column E: only populated for rows 2..3
# Column F: always blank

$data = @(
    @(2,  0, "NSE:CENTENKA",   "NSE:3IINFOLTD",  "", "NSE:ASHOKLEY",   ""),
    @(3,  1, "NSE:DICIND",     "NSE:ACC",        "", "NSE:BALRAMCHIN", ""),
    @(4,  2, "NSE:LAMBODHARA", "NSE:ALMONDZ",    "", "",               ""),
    @(5,  3, "NSE:PAKKA",      "NSE:ARCHIES",    "", "",               ""),
    @(6,  4, "",               "NSE:ASMS",       "", "",               ""),
    @(7,  5, "",               "NSE:BHARATWIRE", "", "",               ""),
    @(8,  6, "",               "NSE:BINANIIND",  "", "",               ""),
    @(9,  7, "",               "NSE:BSOFT",      "", "",               ""),
    @(10, 8, "",               "NSE:CELEBRITY",  "", "",               ""),
    @(11, 9, "",               "NSE:CGCL",       "", "",               ""),
    @(12, 10,"",               "NSE:EDELWEISS",  "", "",               ""),
    @(13, 11,"",               "NSE:EICHERMOT",  "", "",               ""),
    @(14, 12,"",               "NSE:FOODSIN",    "", "",               ""),
    @(15, 13,"",               "NSE:GODREJIND",  "", "",               ""),
    @(16, 14,"",               "NSE:IIFLSEC",    "", "",               ""),
    @(17, 15,"",               "NSE:INDSWFTLTD", "", "",               ""),
    @(18, 16,"",               "NSE:ITDC",       "", "",               ""),
    @(19, 17,"",               "NSE:JUBLINDS",   "", "",               ""),
    @(20, 18,"",               "NSE:MAGADSUGAR", "", "",               ""),
    @(21, 19,"",               "NSE:MANORG",     "", "",               ""),
    @(22, 20,"",               "NSE:MARUTI",     "", "",               ""),
    @(23, 21,"",               "NSE:MATRIMONY",  "", "",               ""),
    @(24, 22,"",               "NSE:MAYURUNIQ",  "", "",               ""),
    @(25, 23,"",               "NSE:METROBRAND", "", "",               ""),
    @(26, 24,"",               "NSE:MOLDTKPAC",  "", "",               ""),
    @(27, 25,"",               "NSE:NOCIL",      "", "",               ""),
    @(28, 26,"",               "NSE:ORIENTCER",  "", "",               ""),
    @(29, 27,"",               "NSE:RANEENGINE", "", "",               ""),
    @(30, 28,"",               "NSE:REMSONSIND", "", "",               ""),
    @(31, 29,"",               "NSE:RHIM",       "", "",               "")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
}
